# Update cryptos list with the latest price/volume snapshot (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain decimal number need an explicit text format while
# being written, otherwise Excel auto-converts the text into a floating point number.
# We apply "@" (Text) format only for the write, then clear the formatting again so the
# cell keeps its original (default/general) style - only the stored value changes.
$numericLookingCells = @{
    'D5' = '609.16'
    'D6' = '147.30'
    'D9' = '0.491'
    'D12' = '0.416'
    'D15' = '30.07'
    'D19' = '11.42'
    'D21' = '15.10'
    'D22' = '432.44'
    'D23' = '0.625'
    'D24' = '79.07'
    'D28' = '8.19'
    'D33' = '25.55'
    'D36' = '7.87'
    'D39' = '5.63'
    'D40' = '173.65'
    'D41' = '0.0856'
    'D42' = '5.24'
    'D45' = '46.12'
    'D49' = '25.05'
    'D50' = '24.00'
    'D51' = '7.22'
}
foreach ($cellRef in $numericLookingCells.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $numericLookingCells[$cellRef]
    $range.ClearFormats()
}

# Remaining cells (plain text, URLs, percentages, and already-non-numeric-looking
# price strings such as "66.611.88") can be written directly.
$ws.Range("D2").Value = '66.611.88'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '3.586.07'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +1.68%  '
$ws.Range("D7").Value = '3.584.69'
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '4.196.24'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '3.590.97'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '66.692.86'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  +2.58%  '
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '3.730.06'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '3.582.18'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("E34").Value = '  -2.48%  '
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  +7.05%  '
$ws.Range("E48").Value = '  -2.79%  '
$ws.Range("E49").Value = '  -3.73%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E51").Value = '  +0.89%  '
